# Corrected index problem with excel.
# Rows 8-17 (participant "A27_19" was missing), causing everything below it
# to be off-by-one in columns A (Participants), E (Free Lies), J (Free Lies Sec.),
# and O (# Free Lies). This shifts those four columns down by one row for rows
# 8-15, inserts the missing participant "F27_19" into row 8 (whose E/J/O data
# was not yet collected), and moves the data that used to overflow past the
# bottom of the table (old row 15, "M31_6") onto row 17, which already holds
# the same participant's other stats but previously lacked this trailing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Participants) -------------------------------------------------
$ws.Range("A8").Value  = "F27_19"
$ws.Range("A9").Value  = "F29_14"
$ws.Range("A10").Value = "M24_18"
$ws.Range("A11").Value = "M25_5"
$ws.Range("A12").Value = "M26_10"
$ws.Range("A13").Value = "M26_12"
$ws.Range("A14").Value = "M26_17"
$ws.Range("A15").Value = "M28_7"
# A17 already is "M31_6" - unchanged

# --- Column E (Free Lies) -----------------------------------------------------
$ws.Range("E9").Value  = -0.1054721431556653
$ws.Range("E10").Value = -0.2480356728751577
$ws.Range("E11").Value = 0.3053562308869803
$ws.Range("E12").Value = 0.07483201451396219
$ws.Range("E13").Value = 0.06230324688604808
$ws.Range("E14").Value = 0.1836924325818537
$ws.Range("E15").Value = 0.1741677061458269
$ws.Range("E17").Value = -0.004623178710214178
$ws.Range("E8").ClearContents()

# --- Column J (Free Lies Sec.) ------------------------------------------------
$ws.Range("J9").Value  = 6.99686465
$ws.Range("J10").Value = 2.3608105
$ws.Range("J11").Value = 2.1692124
$ws.Range("J12").Value = 2.630323366666667
$ws.Range("J13").Value = 3.698001133333334
$ws.Range("J14").Value = 4.58430496
$ws.Range("J15").Value = 7.16894535
$ws.Range("J17").Value = 4.011589625
$ws.Range("J8").ClearContents()

# --- Column O (# Free Lies) ---------------------------------------------------
$ws.Range("O10").Value = 2
$ws.Range("O11").Value = 1
$ws.Range("O12").Value = 6
$ws.Range("O13").Value = 3
$ws.Range("O14").Value = 5
$ws.Range("O15").Value = 4
$ws.Range("O17").Value = 4
$ws.Range("O8").ClearContents()
# O9 stays 2 (unchanged)
